$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the GPS-like coordinate strings in column B (sensor1 / sensor2 rows)
$ws.Range("B3").Value = "45,032N65,06W"
$ws.Range("B4").Value = "31,032N56,43W"

# Widen columns B, D and E to fit the new content / layout
$ws.Columns("B").ColumnWidth = 22.333333333333332
$ws.Columns("D").ColumnWidth = 15.333333333333334
$ws.Columns("E").ColumnWidth = 18.833333333333332

# Touch E7 (leave it blank/default-styled) and give F7 the same number
# format/style as the other cells in column E, extending the used range
# to row 7 / column F
$ws.Range("E7").Value = 5
$ws.Range("E7").ClearContents()
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").NumberFormat = $ws.Range("E1").NumberFormat

# Move the active selection to B4
[void]$ws.Range("B4").Select()
